# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-19) is re-sorted so that the two
# workers' debt periods are grouped in chronological order (1803 before
# 1804) instead of by worker. Only the data cells change - the row
# formatting/styles stay where they are.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (target) contents of the data rows, columns C:G
# (B is always "CC" and does not change)
$ws.Range("C16").Value = "73578626"
$ws.Range("D16").Value = "JOSE DANIEL IGLESIAS ROBLES"
$ws.Range("E16").Value = "1803"
$ws.Range("F16").Value = 48000
$ws.Range("G16").Value = 1200000

$ws.Range("C17").Value = "11077005"
$ws.Range("D17").Value = "RIVER RAMON RUIZ BRAVO"
$ws.Range("E17").Value = "1803"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 781242

$ws.Range("C18").Value = "73578626"
$ws.Range("D18").Value = "JOSE DANIEL IGLESIAS ROBLES"
$ws.Range("E18").Value = "1804"
$ws.Range("F18").Value = 48000
$ws.Range("G18").Value = 1200000

$ws.Range("C19").Value = "11077005"
$ws.Range("D19").Value = "RIVER RAMON RUIZ BRAVO"
$ws.Range("E19").Value = "1804"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 781242
